# Update the simulated transition-probability matrix on Sheet1.
# Source data regenerated after adding more simulated games and
# speeding up the simulate-game logic (see commit message); the
# recomputed probabilities below replace the previous ones cell-by-cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2423076923076923
$ws.Range("C2").Value = 0.4576923076923077
$ws.Range("J2").Value = 0.01923076923076923
$ws.Range("P2").Value = 0.1769230769230769
$ws.Range("S2").Value = 0.1038461538461539
$ws.Range("J3").Value = 0.04958677685950413
$ws.Range("P3").Value = 0.7024793388429752
$ws.Range("S3").Value = 0.2479338842975207
$ws.Range("J4").Value = 0.03846153846153846
$ws.Range("P4").Value = 0.6538461538461539
$ws.Range("S4").Value = 0.3076923076923077
$ws.Range("B6").Value = 0.08411214953271028
$ws.Range("D6").Value = 0.009345794392523364
$ws.Range("F6").Value = 0.102803738317757
$ws.Range("J6").Value = 0.1635514018691589
$ws.Range("O6").Value = 0.01401869158878505
$ws.Range("Q6").Value = 0.2196261682242991
$ws.Range("R6").Value = 0.0514018691588785
$ws.Range("S6").Value = 0.3551401869158878
$ws.Range("B7").Value = 0.0963855421686747
$ws.Range("D7").Value = 0.04216867469879518
$ws.Range("F7").Value = 0.1024096385542169
$ws.Range("J7").Value = 0.1144578313253012
$ws.Range("Q7").Value = 0.1566265060240964
$ws.Range("R7").Value = 0.1325301204819277
$ws.Range("S7").Value = 0.3554216867469879
$ws.Range("B8").Value = 0.07304785894206549
$ws.Range("D8").Value = 0.01007556675062972
$ws.Range("E8").Value = 0.002518891687657431
$ws.Range("F8").Value = 0.04785894206549118
$ws.Range("J8").Value = 0.1209068010075567
$ws.Range("O8").Value = 0.01259445843828715
$ws.Range("Q8").Value = 0.1511335012594459
$ws.Range("R8").Value = 0.1360201511335013
$ws.Range("S8").Value = 0.4458438287153653
$ws.Range("B9").Value = 0.09259259259259259
$ws.Range("D9").Value = 0.01851851851851852
$ws.Range("E9").Value = 0.003703703703703704
$ws.Range("F9").Value = 0.06666666666666667
$ws.Range("J9").Value = 0.1222222222222222
$ws.Range("O9").Value = 0.01851851851851852
$ws.Range("Q9").Value = 0.1592592592592593
$ws.Range("R9").Value = 0.09259259259259259
$ws.Range("S9").Value = 0.4259259259259259
$ws.Range("B10").Value = 0.0947002606429192
$ws.Range("D10").Value = 0.02780191138140747
$ws.Range("F10").Value = 0.07037358818418767
$ws.Range("J10").Value = 0.1181581233709817
$ws.Range("O10").Value = 0.01216333622936577
$ws.Range("Q10").Value = 0.1798436142484796
$ws.Range("R10").Value = 0.0999131190269331
$ws.Range("S10").Value = 0.3970460469157255
$ws.Range("G11").Value = 0.14
$ws.Range("J11").Value = 0.096
$ws.Range("K11").Value = 0.184
$ws.Range("L11").Value = 0.5679999999999999
$ws.Range("S11").Value = 0.012
$ws.Range("G12").Value = 0.8287671232876712
$ws.Range("J12").Value = 0.0958904109589041
$ws.Range("K12").Value = 0.00684931506849315
$ws.Range("L12").Value = 0.03424657534246575
$ws.Range("S12").Value = 0.03424657534246575
$ws.Range("G13").Value = 0.5806451612903226
$ws.Range("J13").Value = 0.2903225806451613
$ws.Range("S13").Value = 0.1290322580645161
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.02
$ws.Range("H15").Value = 0.195
$ws.Range("I15").Value = 0.065
$ws.Range("J15").Value = 0.38
$ws.Range("K15").Value = 0.065
$ws.Range("M15").Value = 0.015
$ws.Range("O15").Value = 0.1
$ws.Range("S15").Value = 0.16
$ws.Range("F16").Value = 0.006329113924050633
$ws.Range("H16").Value = 0.1518987341772152
$ws.Range("I16").Value = 0.1708860759493671
$ws.Range("J16").Value = 0.3670886075949367
$ws.Range("K16").Value = 0.1139240506329114
$ws.Range("M16").Value = 0.01265822784810127
$ws.Range("N16").Value = 0.006329113924050633
$ws.Range("O16").Value = 0.04430379746835443
$ws.Range("S16").Value = 0.1265822784810127
$ws.Range("F17").Value = 0.01574803149606299
$ws.Range("H17").Value = 0.1706036745406824
$ws.Range("I17").Value = 0.1601049868766404
$ws.Range("J17").Value = 0.4120734908136483
$ws.Range("K17").Value = 0.09711286089238845
$ws.Range("M17").Value = 0.007874015748031496
$ws.Range("N17").Value = 0.002624671916010499
$ws.Range("O17").Value = 0.07349081364829396
$ws.Range("S17").Value = 0.06036745406824147
$ws.Range("F18").Value = 0.02212389380530973
$ws.Range("H18").Value = 0.1460176991150443
$ws.Range("I18").Value = 0.1194690265486726
$ws.Range("J18").Value = 0.415929203539823
$ws.Range("K18").Value = 0.07079646017699115
$ws.Range("M18").Value = 0.01327433628318584
$ws.Range("N18").Value = 0.01327433628318584
$ws.Range("O18").Value = 0.0752212389380531
$ws.Range("S18").Value = 0.1238938053097345
$ws.Range("F19").Value = 0.01263689974726201
$ws.Range("H19").Value = 0.2030328559393429
$ws.Range("I19").Value = 0.1213142375737152
$ws.Range("J19").Value = 0.3706823925863522
$ws.Range("K19").Value = 0.09856781802864364
$ws.Range("M19").Value = 0.01684919966301601
$ws.Range("N19").Value = 0.001684919966301601
$ws.Range("O19").Value = 0.06571187868576242
$ws.Range("S19").Value = 0.109519797809604
